$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row of submitted-form data as row 5.
#
# Most values are plain text and can be assigned directly. Column A's value
# is an empty string (not a blank cell) and column C's value is the
# digit-string "20" - both need to be forced to Text so Excel's normal
# "does this look like a number?" auto-detection doesn't turn them into a
# blank/numeric cell (a leading apostrophe is the standard way to tell Excel
# "treat this entry as text").
$ws.Cells.Item(5, 1).Value = "'"
$ws.Cells.Item(5, 2).Value = "احمد"
$ws.Cells.Item(5, 3).Value = "'20"
$ws.Cells.Item(5, 4).Value = "الصمود"
$ws.Cells.Item(5, 5).Value = "الرحلة 1"
$ws.Cells.Item(5, 6).Value = "C1"
$ws.Cells.Item(5, 7).Value = "UNICEF"
$ws.Cells.Item(5, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٤٢:٥٠ م"
